# Update the date line and the 20 answer cells of the division worksheet.
# Each "find" string below is unique within the document at the moment it
# runs, and replacements are ordered so that a freshly-introduced value
# (e.g. "51÷4=12, 3" reused at the end of the table) is never matched
# before the earlier cell holding the same text has already been changed.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Title / date line
Replace-Text "2024-11-06 Wednesday" "2024-11-07 Thursday"

# Row 1
Replace-Text "31÷4=7, 3"  "81÷8=10, 1"
Replace-Text "53÷8=6, 5"  "83÷2=41, 1"
Replace-Text "18÷9=2, 0"  "46÷4=11, 2"
Replace-Text "20÷3=6, 2"  "19÷7=2, 5"
Replace-Text "49÷2=24, 1" "21÷7=3, 0"

# Row 2
Replace-Text "49÷3=16, 1" "21÷4=5, 1"
Replace-Text "65÷8=8, 1"  "66÷5=13, 1"
Replace-Text "50÷2=25, 0" "89÷2=44, 1"
Replace-Text "11÷9=1, 2"  "85÷8=10, 5"
Replace-Text "67÷5=13, 2" "45÷3=15, 0"

# Row 3
Replace-Text "80÷4=20, 0" "70÷9=7, 7"
Replace-Text "29÷8=3, 5"  "89÷5=17, 4"
Replace-Text "95÷9=10, 5" "52÷2=26, 0"
Replace-Text "54÷3=18, 0" "88÷5=17, 3"
Replace-Text "93÷5=18, 3" "87÷9=9, 6"

# Row 4
Replace-Text "95÷2=47, 1" "36÷5=7, 1"
Replace-Text "70÷7=10, 0" "76÷7=10, 6"
Replace-Text "63÷2=31, 1" "15÷3=5, 0"
Replace-Text "27÷2=13, 1" "89÷5=17, 4"
Replace-Text "99÷2=49, 1" "65÷3=21, 2"

# Row 5 (order matters: cell 3's old text "51÷4=12, 3" must be replaced
# before cell 5 is set to that same string)
Replace-Text "35÷5=7, 0"  "58÷4=14, 2"
Replace-Text "19÷3=6, 1"  "58÷4=14, 2"
Replace-Text "51÷4=12, 3" "65÷5=13, 0"
Replace-Text "30÷6=5, 0"  "29÷4=7, 1"
Replace-Text "99÷9=11, 0" "51÷4=12, 3"
